$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column A to fit the new, longer field names (pipeline naming convention)
$ws.Columns.Item(1).ColumnWidth = 46.7109375

# Re-label the pipeline field names in column A.
# The assignment order below matches the order the values were typed into the
# workbook (this controls the order new entries are appended to the shared
# string table).
$ws.Cells.Item(2, 1).Value = "sem_paper_ID"
$ws.Cells.Item(28, 1).Value = "sem_s_STUDY"
$ws.Cells.Item(3, 1).Value = "sem_url"
$ws.Cells.Item(4, 1).Value = "sem_title"
$ws.Cells.Item(5, 1).Value = "sem_venue"
$ws.Cells.Item(11, 1).Value = "sem_year"
$ws.Cells.Item(39, 1).Value = "sem_references[]"
$ws.Cells.Item(6, 1).Value = "sem_publication_venue_name"
$ws.Cells.Item(9, 1).Value = "sem_publication_venue_url"
$ws.Cells.Item(10, 1).Value = "sem_publication_venue_id"
$ws.Cells.Item(35, 1).Value = "sem_external_ids_doi"
$ws.Cells.Item(7, 1).Value = "sem_publication_venue_type"
$ws.Cells.Item(8, 1).Value = "sem_publication_venue_alternate_names"
$ws.Cells.Item(12, 1).Value = "authors_name"
$ws.Cells.Item(13, 1).Value = "authors_sequence"
$ws.Cells.Item(14, 1).Value = "sem_authors_url"
$ws.Cells.Item(17, 1).Value = "sem_authors_homepage"
$ws.Cells.Item(18, 1).Value = "sem_authors_paperCount"
$ws.Cells.Item(19, 1).Value = "sem_authors_citationCount"
$ws.Cells.Item(20, 1).Value = "sem_authors_hIndex"
$ws.Cells.Item(15, 1).Value = "sem_authors_aliases"
$ws.Cells.Item(16, 1).Value = "sem_authors_affiliations"
$ws.Cells.Item(22, 1).Value = "reference_count"
$ws.Cells.Item(23, 1).Value = "citation_count"
$ws.Cells.Item(24, 1).Value = "sem_influential_citation_count"
$ws.Cells.Item(25, 1).Value = "sem_is_openaccess"
$ws.Cells.Item(26, 1).Value = "sem_openaccess_location"
$ws.Cells.Item(27, 1).Value = "General_category"
$ws.Cells.Item(30, 1).Value = "publication_Date"
$ws.Cells.Item(31, 1).Value = "sem_journal_name"
$ws.Cells.Item(32, 1).Value = "sem_journal_volume"
$ws.Cells.Item(34, 1).Value = "sem_external_arxi_id"
$ws.Cells.Item(36, 1).Value = "license_start_date_time"
$ws.Cells.Item(37, 1).Value = "license_start_content_version"
$ws.Cells.Item(38, 1).Value = "license_start_delay_days"
$ws.Cells.Item(33, 1).Value = "sem_citations[]"
$ws.Cells.Item(40, 1).Value = "cross_lang"
$ws.Cells.Item(41, 1).Value = "cross_score"
$ws.Cells.Item(44, 1).Value = "cross_subject"
$ws.Cells.Item(42, 1).Value = "cross_paper_url"
$ws.Cells.Item(46, 1).Value = "cross_issn_type"
$ws.Cells.Item(45, 1).Value = "cross_issn_number"
$ws.Cells.Item(43, 1).Value = "cross_url"

# These two rows reuse already-existing shared strings.
$ws.Cells.Item(21, 1).Value = "abstract"
$ws.Cells.Item(29, 1).Value = "type"

# Update the active selection, as left by the author after editing.
$ws.Range("A49").Select()
